$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from H1 to I1/J1 and set values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Data for columns I (I0) and J (IF), rows 2-13
$data = @(
    @(8, 9),
    @(9, 9),
    @(7, 9),
    @(5, 8),
    @(7, 8),
    @(2, 5),
    @(4, 7),
    @(8, 8),
    @(8, 9),
    @(7, 9),
    @(9, 9),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
